# A new September transaction was recorded ("balance your axis" at
# 2024-09-25 11:06:04), which is more recent than the existing most-recent
# entry in row 48. Insert a new row above row 48 on the "2024" sheet,
# pushing all subsequent rows (old row 48 ... old row 210) down by one,
# then populate the new row's September_Details / September_Date cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Rows(48).Insert()

$ws.Range("R48").Value = "balance your axis"
$ws.Range("S48").Value = "2024-09-25 11:06:04"
